$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# Row 2
Set-TextValue $ws.Range("D2") '70.059.82'
Set-TextValue $ws.Range("E2") '  +3.16%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.455.16'
Set-TextValue $ws.Range("E3") '  +1.50%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.06%  '

# Row 5
Set-TextValue $ws.Range("D5") '567.71'
Set-TextValue $ws.Range("E5") '  +2.17%  '

# Row 6
Set-TextValue $ws.Range("D6") '167.64'
Set-TextValue $ws.Range("E6") '  +4.37%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.512'
Set-TextValue $ws.Range("E8") '  +0.25%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.176'
Set-TextValue $ws.Range("E9") '  +13.13%  '

# Row 10
Set-TextValue $ws.Range("D10") '2.451.87'
Set-TextValue $ws.Range("E10") '  +1.32%  '

# Row 11
Set-TextValue $ws.Range("E11") '  -1.55%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +3.28%  '

# Row 13
Set-TextValue $ws.Range("D13") '4.71'
Set-TextValue $ws.Range("E13") '  -0.89%  '

# Row 14
Set-TextValue $ws.Range("E14") '  +8.97%  '

# Row 15
Set-TextValue $ws.Range("D15") '69.942.52'
Set-TextValue $ws.Range("E15") '  +3.07%  '

# Row 16
Set-TextValue $ws.Range("D16") '2.906.37'
Set-TextValue $ws.Range("E16") '  +0.02%  '

# Row 17
Set-TextValue $ws.Range("D17") '24.17'
Set-TextValue $ws.Range("E17") '  +5.31%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.458.79'
Set-TextValue $ws.Range("E18") '  +0.45%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.88'
Set-TextValue $ws.Range("E19") '  +5.86%  '

# Row 20
Set-TextValue $ws.Range("D20") '7.21'
Set-TextValue $ws.Range("E20") '  +6.15%  '

# Row 21
Set-TextValue $ws.Range("D21") '343.63'
Set-TextValue $ws.Range("E21") '  +2.69%  '

# Row 22
Set-TextValue $ws.Range("E22") '  +3.78%  '

# Row 23
$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D23") '2.01'
Set-TextValue $ws.Range("E23") '  +8.23%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D24") '1.00'
Set-TextValue $ws.Range("E24") '  -0.02%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D25") '66.52'
Set-TextValue $ws.Range("E25") '  -0.14%  '

# Row 26
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D26") '3.89'
Set-TextValue $ws.Range("E26") '  +7.83%  '

# Row 27
$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws.Range("D27") '2.582.94'
Set-TextValue $ws.Range("E27") '  +1.49%  '

# Row 28
Set-TextValue $ws.Range("E28") '  +6.64%  '

# Row 29
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D29") '1.00'
Set-TextValue $ws.Range("E29") '  +0.14%  '

# Row 30
Set-TextValue $ws.Range("E30") '  +8.78%  '

# Row 31
Set-TextValue $ws.Range("E31") '  +4.40%  '

# Row 32
Set-TextValue $ws.Range("E32") '  +11.58%  '

# Row 33
Set-TextValue $ws.Range("D33") '451.70'
Set-TextValue $ws.Range("E33") '  +7.67%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.999'
Set-TextValue $ws.Range("E34") '  -0.17%  '

# Row 35
Set-TextValue $ws.Range("E35") '  +2.62%  '

# Row 36
Set-TextValue $ws.Range("D36") '161.46'
Set-TextValue $ws.Range("E36") '  +0.39%  '

# Row 37
Set-TextValue $ws.Range("D37") '19.09'
Set-TextValue $ws.Range("E37") '  +0.74%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.111'
Set-TextValue $ws.Range("E38") '  +7.36%  '

# Row 39
Set-TextValue $ws.Range("E39") '  +0.03%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +3.03%  '

# Row 41
Set-TextValue $ws.Range("E41") '  +4.82%  '

# Row 42
Set-TextValue $ws.Range("E42") '  +6.96%  '

# Row 43
Set-TextValue $ws.Range("D43") '4.45'
Set-TextValue $ws.Range("E43") '  +4.46%  '

# Row 44
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D44") '37.80'
Set-TextValue $ws.Range("E44") '  +0.59%  '

# Row 45
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D45") '1.09'
Set-TextValue $ws.Range("E45") '  +5.21%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D46") '2.15'
Set-TextValue $ws.Range("E46") '  +8.60%  '

# Row 47
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D47") '3.41'
Set-TextValue $ws.Range("E47") '  +2.76%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D48") '133.51'
Set-TextValue $ws.Range("E48") '  +4.39%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D49") '0.0725'
Set-TextValue $ws.Range("E49") '  +1.91%  '

# Row 50
$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D50") '0.492'
Set-TextValue $ws.Range("E50") '  +4.19%  '

# Row 51
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D51") '0.564'
Set-TextValue $ws.Range("E51") '  +1.99%  '

Write-Host "Applied cryptos update"